$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of row 14 into row 15 so date/number formatting match
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A15").Value = 45600
$ws.Range("B15").Value = 2.15
$ws.Range("C15").Value = "Besprechung"
$ws.Range("D15").Value = "Aufgabenverteilung, Organisierung "

$ws.Range("F22").Select()
